# ApplitoolKit Dashboard update
# 1. Insert a new "launchUrl" worksheet as the first sheet, containing a
#    small two-row/two-column table with a header row and a hyperlinked URL.
# 2. Clear the "Add" sheet's former "selected" state / stale active-cell
#    selection now that "launchUrl" is the tab shown on open.

$wb = $excel.ActiveWorkbook

# --- Create the new "launchUrl" worksheet --------------------------------
$launchSheet = $wb.Worksheets.Add()
$launchSheet.Name = "launchUrl"
$launchSheet.Move($wb.Worksheets.Item(1))

# NOTE: fetch sheet references used below only *after* the Move() call
# above, since moving a sheet invalidates references obtained earlier.
$addSheet = $wb.Worksheets.Item("Add")

# Data
$launchSheet.Range("A1").Value = "_rowId"
$launchSheet.Range("B1").Value = "Navigate@URL"
$launchSheet.Range("A2").Value = "_DefaultCompRow"
$launchSheet.Range("B2").Value = "http://juliemr.github.io/protractor-demo/"

# Reuse the same bold header style already used on the other sheets by
# copying the formatting from the "Add" header cell.
$addSheet.Range("A1").Copy()
$launchSheet.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Update the previously-first sheet's view state ----------------------
# The "Add" sheet used to be the active/selected tab with A3 selected;
# now that launchUrl takes over as the opening tab, select the header row
# on "Add" instead (matches the post-edit workbook) and let tab-selection
# move naturally once "launchUrl" is activated below.
$addSheet.Activate()
$addSheet.Range("A1:B1").Select()

# Hyperlink the URL cell.
$launchSheet.Hyperlinks.Add($launchSheet.Range("B2"), "http://juliemr.github.io/protractor-demo/")

# Column widths sized to fit the content (header / url text).
$launchSheet.Columns.Item(1).ColumnWidth = 15.8307291667
$launchSheet.Columns.Item(2).ColumnWidth = 35.1666666667

# Freeze the header row and leave the active cell parked at D5, as on the
# other sheets.
$launchSheet.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$launchSheet.Range("D5").Select()

# Make the new sheet the active / selected tab.
$launchSheet.Activate()
